# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets, row 3 (F3) and row 5 (F5), to match the newly scraped data.
$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 402
    $ws.Range("F5").Value = 111
}
